# luban: remove the --export_test_data option, add the --export_exclude_tags
# option ("delete --export_test_data option, add --export_exclude_tags
# option").
#
# In the tb_role_csv tag sample sheet this means:
#   - the two obsolete "faLse" / "false" tag sample rows are removed
#   - the remaining rows shift up to fill the gap
#   - the last sample row (previously "测试" / "测试") is repurposed into a
#     generic "any" / "any" example row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete rows (the "faLse" and "false" tag samples).
# This shifts every following row up by two automatically.
$null = $ws.Rows("8:9").Delete()

# The row that used to hold the "测试" / "测试" example (now row 12 after the
# shift) becomes a generic "any" / "any" example row.
$ws.Range("A12").Value = "any"
$ws.Range("C12").Value = "any"

# Move the active selection to A9, matching the post-edit cursor position.
$null = $ws.Range("A9").Select()
